# Auto-generated edit script: Add data for 2024-08-11
# Updates 2024 year-to-date (column K) totals, plus a couple of late corrections
# to prior columns (B, J), across the Citywide Totals, By Neighborhood, and per-
# neighborhood sheets, matching the upstream dataset refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 4873
$ws.Range("J3").Value = 8078
$ws.Range("K3").Value = 5018
$ws.Range("B4").Value = 1703
$ws.Range("J4").Value = 1832
$ws.Range("K4").Value = 1040
$ws.Range("K6").Value = 5621
$ws.Range("B7").Value = 23336
$ws.Range("K7").Value = 16908

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K7").Value = 499
$ws.Range("K8").Value = 1130
$ws.Range("K11").Value = 331
$ws.Range("K13").Value = 17
$ws.Range("K15").Value = 168
$ws.Range("K19").Value = 512
$ws.Range("K20").Value = 390
$ws.Range("K23").Value = 170
$ws.Range("K25").Value = 81
$ws.Range("K29").Value = 908
$ws.Range("K30").Value = 65
$ws.Range("K31").Value = 186
$ws.Range("K33").Value = 716
$ws.Range("K36").Value = 222
$ws.Range("K37").Value = 571
$ws.Range("K38").Value = 16
$ws.Range("K40").Value = 40
$ws.Range("K42").Value = 629
$ws.Range("K43").Value = 150
$ws.Range("K44").Value = 148
$ws.Range("K47").Value = 116
$ws.Range("K49").Value = 95
$ws.Range("K51").Value = 214
$ws.Range("K52").Value = 441
$ws.Range("K55").Value = 195
$ws.Range("B63").Value = 407
$ws.Range("K63").Value = 49
$ws.Range("K64").Value = 107
$ws.Range("K65").Value = 382
$ws.Range("K67").Value = 649
$ws.Range("K73").Value = 142
$ws.Range("K75").Value = 58
$ws.Range("K76").Value = 231
$ws.Range("K78").Value = 197
$ws.Range("K79").Value = 414
$ws.Range("K83").Value = 371
$ws.Range("K84").Value = 126
$ws.Range("K85").Value = 772
$ws.Range("K89").Value = 243
$ws.Range("K91").Value = 183
$ws.Range("K92").Value = 62
$ws.Range("K93").Value = 63
$ws.Range("K94").Value = 222
$ws.Range("K97").Value = 135
$ws.Range("K98").Value = 83
$ws.Range("K100").Value = 32
$ws.Range("B101").Value = 23336
$ws.Range("K101").Value = 16908

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K3").Value = 159
$ws.Range("K6").Value = 129
$ws.Range("K7").Value = 499

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K6").Value = 113
$ws.Range("K7").Value = 331

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K2").Value = 65
$ws.Range("K3").Value = 75
$ws.Range("K6").Value = 76
$ws.Range("K7").Value = 243

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K3").Value = 262
$ws.Range("K6").Value = 183
$ws.Range("K7").Value = 772

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K3").Value = 118
$ws.Range("K4").Value = 25
$ws.Range("K7").Value = 441

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 317
$ws.Range("K6").Value = 382
$ws.Range("K7").Value = 1130

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K3").Value = 135
$ws.Range("K6").Value = 86
$ws.Range("K7").Value = 371

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 197
$ws.Range("K3").Value = 267
$ws.Range("K6").Value = 205
$ws.Range("K7").Value = 716

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 162
$ws.Range("K3").Value = 189
$ws.Range("K7").Value = 571

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K3").Value = 97
$ws.Range("K7").Value = 382

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("K2").Value = 18
$ws.Range("K7").Value = 65

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K2").Value = 61
$ws.Range("K3").Value = 43
$ws.Range("K7").Value = 186

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 186
$ws.Range("K3").Value = 227
$ws.Range("K4").Value = 37
$ws.Range("K7").Value = 649

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K2").Value = 38
$ws.Range("K7").Value = 126

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K4").Value = 9
$ws.Range("K7").Value = 95

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 263
$ws.Range("K3").Value = 325
$ws.Range("K6").Value = 251
$ws.Range("K7").Value = 908

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 162
$ws.Range("K7").Value = 512

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K2").Value = 34
$ws.Range("K7").Value = 148

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K4").Value = 15
$ws.Range("K7").Value = 231

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K3").Value = 199
$ws.Range("K7").Value = 629

$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("K2").Value = 2
$ws.Range("K6").Value = 17

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K6").Value = 71
$ws.Range("K7").Value = 197

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K2").Value = 57
$ws.Range("K3").Value = 57
$ws.Range("K7").Value = 195

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K4").Value = 10
$ws.Range("K7").Value = 170

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K3").Value = 86
$ws.Range("K7").Value = 183

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 137
$ws.Range("K3").Value = 132
$ws.Range("K7").Value = 414

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K3").Value = 31
$ws.Range("K6").Value = 40
$ws.Range("K7").Value = 107

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 131
$ws.Range("K6").Value = 111
$ws.Range("K7").Value = 390

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 86
$ws.Range("K3").Value = 63
$ws.Range("K6").Value = 52
$ws.Range("K7").Value = 222

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("K2").Value = 23
$ws.Range("K7").Value = 63

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("K6").Value = 19
$ws.Range("K7").Value = 32

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K3").Value = 40
$ws.Range("K6").Value = 93
$ws.Range("K7").Value = 222

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("K2").Value = 26
$ws.Range("K7").Value = 81

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K2").Value = 33
$ws.Range("K7").Value = 116

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K2").Value = 58
$ws.Range("K7").Value = 168

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K2").Value = 13
$ws.Range("K7").Value = 83

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J3").Value = 74
$ws.Range("J4").Value = 18
$ws.Range("K4").Value = 10
$ws.Range("K7").Value = 142

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K6").Value = 80
$ws.Range("K7").Value = 135

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("K3").Value = 11
$ws.Range("K7").Value = 62

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("K6").Value = 9
$ws.Range("K7").Value = 58

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K6").Value = 73
$ws.Range("K7").Value = 214

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K2").Value = 29
$ws.Range("K4").Value = 19
$ws.Range("K7").Value = 150

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("K3").Value = 18
$ws.Range("K7").Value = 40

$ws = $wb.Worksheets.Item("Grant Park")
$ws.Range("K5").Value = 7
$ws.Range("K6").Value = 16
